$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CB")

# Operating Income (row 6) and Non-operating Income/Expense (row 7)
$ws.Range("B6").Value = 6839000000.0
$ws.Range("B7").Value = -506000000.0

# Gross Margin (row 14)
$ws.Range("D14").Value = 0.1706
$ws.Range("E14").Value = 0.1857
$ws.Range("F14").Value = 0.2331
$ws.Range("G14").Value = 0.2505

# EBIT Margin (row 15)
$ws.Range("D15").Value = 0.0943
$ws.Range("E15").Value = 0.0954
$ws.Range("F15").Value = 0.1488
$ws.Range("G15").Value = 0.1704

# EBT margin (row 16)
$ws.Range("D16").Value = 0.0789
$ws.Range("E16").Value = 0.0795
$ws.Range("F16").Value = 0.1322
$ws.Range("G16").Value = 0.1535

# Net Profit Margin (row 17)
$ws.Range("D17").Value = 0.0656
$ws.Range("F17").Value = 0.108
$ws.Range("G17").Value = 0.1303

# Free Cash Flow Margin (row 18)
$ws.Range("D18").Value = 0.2653
$ws.Range("E18").Value = 0.2257
$ws.Range("F18").Value = 0.2045
$ws.Range("G18").Value = 0.1855

# EBITDA (row 19) - was blank inline string, now numeric
$ws.Range("B19").Value = 7484000000.0

# EBIT (row 20)
$ws.Range("B20").Value = 6839000000.0

# EPS (Diluted, from Cont. Ops) (row 21)
$ws.Range("B21").Value = 12.3105

# EPS (Basic, Consolidated) (row 22) - was blank inline string, now numeric
$ws.Range("B22").Value = 12.3705

# EPS (Basic, from Continuous Ops) (row 23) - was blank inline string, now numeric
$ws.Range("B23").Value = 12.3705

# EBITDA Margin (row 26)
$ws.Range("D26").Value = 0.1136
$ws.Range("E26").Value = 0.115
$ws.Range("F26").Value = 0.1684
$ws.Range("G26").Value = 0.1908

# Operating Cash Flow Margin (row 27)
$ws.Range("D27").Value = 0.2653
$ws.Range("E27").Value = 0.2257
$ws.Range("F27").Value = 0.2045
$ws.Range("G27").Value = 0.1855

$wb.Save()
